$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 filters (D2): adjust the PRODUCTS grid filter string
$ws.Range("D2").Value = "i=2&p=25&s%5BL3%5D=-1&f%5Bscope%5D%5Bvalue%5D=PRODUCT_CATALOG&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=1&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=5&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"

# Fix row 10 (UNASSIGNED PRODUCTS AND ITEMS) columns, which incorrectly pointed at another grid's columns
$ws.Range("B10").Value = "identifier,PRODUCT_DESCRIPTION,BASE_UOM,UOM_QTY_FACTOR"

# Add new row 11: Products - COVID (WHO's COVID-19 product list)
$ws.Range("A11").Value = "Products - COVID"
$ws.Range("B11").Value = "WHO_COVID19_ITEM_CODE,PRODUCT_DESCRIPTION,complete_variant_products,price_reference"
$ws.Range("C11").Value = "product-grid"
$ws.Range("D11").Value = "i=1&p=25&s%5BWHO_COVID19_ITEM_CODE%5D=-1&f%5Bscope%5D%5Bvalue%5D=PRODUCT_CATALOG&f%5Bcategory%5D%5Bvalue%5D%5BtreeId%5D=431&f%5Bcategory%5D%5Bvalue%5D%5BcategoryId%5D=0&f%5Bcategory%5D%5Btype%5D=1&t=product-grid"
$ws.Range("E11").Value = "admin"
$ws.Range("F11").Value = "public"
